# "update metric for ep06" -- append the new EP05 data row to the Data
# table, grow the Table1 ListObject to cover it, and move the selection
# the way the author left it when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New row of metrics (Episode, Live, Views, Feedback, Discussions, Stars, Total Views)
$ws.Range("A6").Value = "EP05"
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 14
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 128

# Grow the structured table (and its autofilter) to include the new row.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:G6"))

# Leave the selection where the author left it on save.
[void]$ws.Range("C30").Select()
